$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2..141 (columns A..R): weekly refresh of
# the "Puerro" price series - one new observation added and the existing
# observations re-sequenced / updated accordingly.
$data = @(
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45203, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44208, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44699, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44491, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44160, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 7000, 8000, 7600, '$/paquete 20 unidades', 'Provincia de Chacabuco', 380, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44349, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 130, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45196, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44292, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44252, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44526, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 170, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45070, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44342, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44847, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 11000, 12000, 11571, '$/paquete 20 unidades', 'Provincia de Melipilla', 579, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44847, 13, 100112005, 'Puerro', 'Sin especificar', 'Segunda', 40, 10000, 10000, 10000, '$/paquete 20 unidades', 'Provincia de Melipilla', 500, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44189, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44244, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44519, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44358, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7500, 8000, 7750, '$/paquete 20 unidades', 'Provincia de Chacabuco', 388, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44301, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44295, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44421, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 180, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44463, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7500, 8000, 7750, '$/paquete 20 unidades', 'Provincia de Chacabuco', 388, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44510, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44484, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44720, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44650, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 9000, 10000, 9500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 475, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44266, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45119, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44299, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44442, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 180, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44539, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 133, 6000, 7000, 6504, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44497, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 180, 6000, 7000, 6556, '$/paquete 20 unidades', 'Provincia de Chacabuco', 328, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44433, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 142, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44873, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 9000, 8571, '$/paquete 20 unidades', 'Provincia de Chacabuco', 429, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44321, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 250, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44671, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44419, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44435, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 302, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45168, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44363, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44482, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44708, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 80, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44707, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 6000, 7000, 6571, '$/paquete 20 unidades', 'Provincia de Chacabuco', 329, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45175, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44166, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44267, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44657, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44825, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 12000, 12000, 12000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 600, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44552, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 106, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44286, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44505, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45126, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45049, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 9000, 9000, 9000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 450, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44259, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44643, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 9000, 8500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 425, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45030, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 9000, 9000, 9000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 450, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44447, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 106, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44517, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44468, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 133, 7000, 8000, 7504, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44477, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44162, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44272, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45042, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44308, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45210, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45015, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44664, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 106, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44384, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 9000, 8500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 425, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44615, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 79, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44623, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 106, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45084, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44524, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44573, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 106, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44232, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 60, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44860, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 9000, 10000, 9571, '$/paquete 20 unidades', 'Provincia de Chacabuco', 479, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44692, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 124, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44278, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44679, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 97, 8000, 9000, 8505, '$/paquete 20 unidades', 'Provincia de Chacabuco', 425, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44328, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44314, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44405, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44398, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7500, 8000, 7750, '$/paquete 20 unidades', 'Provincia de Chacabuco', 388, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44391, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 52, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44489, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45154, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44461, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 79, 7000, 8000, 7494, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44428, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 97, 8000, 9000, 8505, '$/paquete 20 unidades', 'Provincia de Chacabuco', 425, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44302, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44832, 13, 100112005, 'Puerro', 'Sin especificar', 'Segunda', 30, 12000, 12000, 12000, '$/paquete 20 unidades', 'Provincia de Melipilla', 600, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44281, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 250, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44203, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 7000, 8000, 7400, '$/paquete 20 unidades', 'Provincia de Chacabuco', 370, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44344, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 210, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44426, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 97, 7000, 8000, 7505, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44755, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44776, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44769, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45077, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45217, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44407, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44356, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44412, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7500, 8000, 7750, '$/paquete 20 unidades', 'Provincia de Chacabuco', 388, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45112, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44229, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44855, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 65, 9000, 10000, 9538, '$/paquete 20 unidades', 'Provincia de Melipilla', 477, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44855, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 80, 10000, 10000, 10000, '$/paquete 20 unidades', 'Provincia de Santiago', 500, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44365, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 180, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44475, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44748, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44545, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 6000, 7000, 6500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 325, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44195, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44636, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 97, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44306, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44503, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 97, 7000, 8000, 7505, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44685, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44790, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45091, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44215, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 80, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45140, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45147, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44224, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 120, 6000, 7000, 6667, '$/paquete 20 unidades', 'Provincia de Chacabuco', 333, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44762, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44630, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 79, 9000, 10000, 9494, '$/paquete 20 unidades', 'Provincia de Chacabuco', 475, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44273, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44454, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44370, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7500, 8000, 7750, '$/paquete 20 unidades', 'Provincia de Chacabuco', 388, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44239, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44214, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44265, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44335, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 250, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45063, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 7000, 7000, 7000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 350, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45133, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44174, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44218, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 80, 6000, 7000, 6625, '$/paquete 20 unidades', 'Provincia de Chacabuco', 331, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45028, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44329, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44498, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 130, 9000, 10000, 9615, '$/paquete 20 unidades', 'Provincia de Santiago', 481, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45035, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 70, 9000, 9000, 9000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 450, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 45224, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 160, 7000, 8000, 7500, '$/paquete 20 unidades', 'Provincia de Chacabuco', 375, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44167, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 50, 8000, 8000, 8000, '$/paquete 20 unidades', 'Provincia de Chacabuco', 400, 20, 'Hortaliza'),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44414, 13, 100112005, 'Puerro', 'Sin especificar', 'Primera', 180, 7500, 8000, 7750, '$/paquete 20 unidades', 'Provincia de Chacabuco', 388, 20, 'Hortaliza')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $targetRow = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $row[$c]
    }
}

# Make sure the date column (D) keeps the date-time number format used
# throughout the sheet, including on the newly appended row 141.
$ws.Range("D2:D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Updated $($data.Length) data rows (2..141)"
